# Applies the "progress: find bks, comparison table and statistical analysis
# working" update: a new randomized run's results are written into the
# workbook (Resumen / Solucion / Metricas sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Resumen" sheet -> update the Maximo value in C2
# ---------------------------------------------------------------------------
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 627.8886906692617

# ---------------------------------------------------------------------------
# 2) "Solucion" sheet -> rewrite the Pedido/Salida assignment table (A2:B41)
# ---------------------------------------------------------------------------
$wsSolucion = $wb.Worksheets.Item("Solucion")

$data = @(
    @("Pedido_4", "S001"),
    @("Pedido_22", "S025"),
    @("Pedido_11", "S005"),
    @("Pedido_6", "S029"),
    @("Pedido_12", "S002"),
    @("Pedido_5", "S026"),
    @("Pedido_37", "S006"),
    @("Pedido_32", "S003"),
    @("Pedido_14", "S030"),
    @("Pedido_24", "S007"),
    @("Pedido_25", "S027"),
    @("Pedido_38", "S004"),
    @("Pedido_7", "S031"),
    @("Pedido_40", "S008"),
    @("Pedido_19", "S009"),
    @("Pedido_33", "S028"),
    @("Pedido_3", "S013"),
    @("Pedido_28", "S032"),
    @("Pedido_35", "S010"),
    @("Pedido_18", "S033"),
    @("Pedido_1", "S037"),
    @("Pedido_27", "S014"),
    @("Pedido_16", "S034"),
    @("Pedido_29", "S011"),
    @("Pedido_36", "S015"),
    @("Pedido_2", "S012"),
    @("Pedido_31", "S038"),
    @("Pedido_15", "S035"),
    @("Pedido_34", "S016"),
    @("Pedido_17", "S017"),
    @("Pedido_23", "S039"),
    @("Pedido_30", "S021"),
    @("Pedido_13", "S036"),
    @("Pedido_10", "S018"),
    @("Pedido_9", "S040"),
    @("Pedido_8", "S022"),
    @("Pedido_39", "S019"),
    @("Pedido_20", "S023"),
    @("Pedido_21", "S020"),
    @("Pedido_26", "S024")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $data[$i][0]
    $wsSolucion.Cells.Item($row, 2).Value = $data[$i][1]
}

# ---------------------------------------------------------------------------
# 3) "Metricas" sheet -> update the Tiempo values for Z1 and Z2
# ---------------------------------------------------------------------------
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 627.8886906692617
$wsMetricas.Range("B3").Value = 491.7017785706563
